$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5739
$ws1.Range("F3").Value = 873
$ws1.Range("G3").Value = "已售罄"
$ws1.Range("F4").Value = 88
$ws1.Range("F6").Value = 9

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5739
$ws4.Range("F3").Value = 873
$ws4.Range("G3").Value = "已售罄"
$ws4.Range("F4").Value = 88
$ws4.Range("F7").Value = 9
